$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (ownTeam, oppTeam) so the old
# D:I (batsman..sr) shifts right to F:K.
$ws.Range("D1:E1").EntireColumn.Insert() | Out-Null

# Header row
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# All data cells in this sheet are stored as text (numberStoredAsText),
# including the numeric-looking ones, so force text formatting before
# assigning values to keep them as text rather than real numbers.
$ws.Range("G2:K3").NumberFormat = "@"

# Data rows are also reordered: the former row 3 (Sharjah) becomes row 2,
# and the former row 2 (Abu Dhabi) becomes row 3. Write the full final
# grid explicitly to match both the reorder and the new team columns.
$ws.Range("A2").Value = " Sharjah"
$ws.Range("B2").Value = " November 03 2020"
$ws.Range("C2").Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Sunrisers Hyderabad"
$ws.Range("F2").Value = "James Pattinson "
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "5"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "80.00"

$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " September 19 2020"
$ws.Range("C3").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D3").Value = "Mumbai Indians"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "James Pattinson "
$ws.Range("G3").Value = "11"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "2"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "137.50"
